{"js": "// Replace the twenty-five \"NN\u00d7NN=NNNN\" answer strings scattered across the\n// table cells with their new values. Every source string is unique in the\n// document, so a direct search-and-replace per pair is safe and unambiguous.\nconst replacements = [\n  [\"20\u00d776=1520\", \"38\u00d769=2622\"],\n  [\"79\u00d799=7821\", \"65\u00d767=4355\"],\n  [\"51\u00d729=1479\", \"22\u00d798=2156\"],\n  [\"35\u00d721=735\", \"45\u00d713=585\"],\n  [\"86\u00d772=6192\", \"38\u00d730=1140\"],\n  [\"61\u00d782=5002\", \"48\u00d753=2544\"],\n  [\"45\u00d759=2655\", \"23\u00d766=1518\"],\n  [\"32\u00d719=608\", \"80\u00d740=3200\"],\n  [\"92\u00d744=4048\", \"27\u00d770=1890\"],\n  [\"12\u00d714=168\", \"35\u00d723=805\"],\n  [\"28\u00d726=728\", \"54\u00d753=2862\"],\n  [\"43\u00d737=1591\", \"87\u00d761=5307\"],\n  [\"24\u00d717=408\", \"64\u00d798=6272\"],\n  [\"31\u00d799=3069\", \"87\u00d761=5307\"],\n  [\"99\u00d751=5049\", \"35\u00d791=3185\"],\n  [\"51\u00d724=1224\", \"72\u00d727=1944\"],\n  [\"68\u00d773=4964\", \"95\u00d792=8740\"],\n  [\"17\u00d748=816\", \"12\u00d725=300\"],\n  [\"11\u00d747=517\", \"43\u00d711=473\"],\n  [\"18\u00d739=702\", \"20\u00d738=760\"],\n  [\"83\u00d711=913\", \"19\u00d769=1311\"],\n  [\"39\u00d748=1872\", \"23\u00d737=851\"],\n  [\"58\u00d750=2900\", \"52\u00d793=4836\"],\n  [\"86\u00d799=8514\", \"42\u00d772=3024\"],\n  [\"19\u00d768=1292\", \"51\u00d742=2142\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the twenty-five \"NN\u00d7NN=NNNN\" answer strings scattered across the\n# table cells with their new values. Every source string is unique in the\n# document, so a Find/Replace pass per pair is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"20\u00d776=1520\"; New = \"38\u00d769=2622\" },\n    @{ Old = \"79\u00d799=7821\"; New = \"65\u00d767=4355\" },\n    @{ Old = \"51\u00d729=1479\"; New = \"22\u00d798=2156\" },\n    @{ Old = \"35\u00d721=735\";  New = \"45\u00d713=585\"  },\n    @{ Old = \"86\u00d772=6192\"; New = \"38\u00d730=1140\" },\n    @{ Old = \"61\u00d782=5002\"; New = \"48\u00d753=2544\" },\n    @{ Old = \"45\u00d759=2655\"; New = \"23\u00d766=1518\" },\n    @{ Old = \"32\u00d719=608\";  New = \"80\u00d740=3200\" },\n    @{ Old = \"92\u00d744=4048\"; New = \"27\u00d770=1890\" },\n    @{ Old = \"12\u00d714=168\";  New = \"35\u00d723=805\"  },\n    @{ Old = \"28\u00d726=728\";  New = \"54\u00d753=2862\" },\n    @{ Old = \"43\u00d737=1591\"; New = \"87\u00d761=5307\" },\n    @{ Old = \"24\u00d717=408\";  New = \"64\u00d798=6272\" },\n    @{ Old = \"31\u00d799=3069\"; New = \"87\u00d761=5307\" },\n    @{ Old = \"99\u00d751=5049\"; New = \"35\u00d791=3185\" },\n    @{ Old = \"51\u00d724=1224\"; New = \"72\u00d727=1944\" },\n    @{ Old = \"68\u00d773=4964\"; New = \"95\u00d792=8740\" },\n    @{ Old = \"17\u00d748=816\";  New = \"12\u00d725=300\"  },\n    @{ Old = \"11\u00d747=517\";  New = \"43\u00d711=473\"  },\n    @{ Old = \"18\u00d739=702\";  New = \"20\u00d738=760\"  },\n    @{ Old = \"83\u00d711=913\";  New = \"19\u00d769=1311\" },\n    @{ Old = \"39\u00d748=1872\"; New = \"23\u00d737=851\"  },\n    @{ Old = \"58\u00d750=2900\"; New = \"52\u00d793=4836\" },\n    @{ Old = \"86\u00d799=8514\"; New = \"42\u00d772=3024\" },\n    @{ Old = \"19\u00d768=1292\"; New = \"51\u00d742=2142\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $range.Find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
